$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Clasificación" (standings) table with this matchday's results.
# Columns (1-based, A=1): C=3 PJ, D=4 PG, E=5 PE, F=6 PP, G=7 TD, H=8 TP,
#                         I=9 DT, J=10 V1, K=11 V2, L=12 VC, M=13 PTS
$data = @(
    @(2.0, 1.0, 0.0, 1.0, 1.0, 1.0, 0.0, 1.0, 0.0, 0.0, 3.0),  # row 2 - David
    @(2.0, 0.0, 0.0, 2.0, 1.0, 5.0, -4.0, 0.0, 0.0, 0.0, 0.0), # row 3 - Pedro
    @(2.0, 2.0, 0.0, 0.0, 3.0, 0.0, 3.0, 1.0, 1.0, 0.0, 7.0),  # row 4 - Adonay
    @(2.0, 0.0, 1.0, 1.0, 0.0, 1.0, -1.0, 0.0, 0.0, 0.0, 1.0), # row 5 - Richard
    @(2.0, 1.0, 0.0, 1.0, 3.0, 2.0, 1.0, 0.0, 0.0, 1.0, 5.0),  # row 6 - Iván
    @(2.0, 1.0, 1.0, 0.0, 1.0, 0.0, 1.0, 1.0, 0.0, 0.0, 4.0),  # row 7 - Nico
    @(2.0, 1.0, 0.0, 1.0, 1.0, 1.0, 0.0, 1.0, 0.0, 0.0, 3.0),  # row 8 - Nicolás
    @(2.0, 1.0, 0.0, 1.0, 1.0, 1.0, 0.0, 1.0, 0.0, 0.0, 3.0)   # row 9 - Vicente
)

for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $i + 2
    $values = $data[$i]
    for ($j = 0; $j -lt $values.Length; $j++) {
        $col = $j + 3
        $ws.Cells.Item($row, $col).Value = $values[$j]
    }
}
